# Updated cryptos list on Thu Feb 15 15:05:31 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row
# with the latest scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell as literal TEXT
# (matches the source data, e.g. '1.00' / '3.30' / '70.79') without
# leaving the cell's number format changed.
function Set-TextValue($cell, $text) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Formula = $text
    $range.Style = "Normal"
}

# row => @{ Price = ...; Volume = ... } (Price omitted where unchanged)
$updates = [ordered]@{
    2 = @{ Price = "52.616.32"; Volume = "  +1.44%  " }
    3 = @{ Price = "2.853.37"; Volume = "  +3.47%  " }
    4 = @{ Price = "1.00"; Volume = "  +0.03%  " }
    5 = @{ Price = "363.21"; Volume = "  +9.30%  " }
    6 = @{ Price = "116.88"; Volume = "  -0.81%  " }
    7 = @{ Price = "0.553"; Volume = "  +3.62%  " }
    8 = @{ Volume = "  -0.05%  " }
    9 = @{ Volume = "  +5.11%  " }
    10 = @{ Price = "43.05"; Volume = "  +2.88%  " }
    11 = @{ Volume = "  +4.39%  " }
    12 = @{ Volume = "  +0.63%  " }
    13 = @{ Volume = "  +1.55%  " }
    14 = @{ Price = "7.92"; Volume = "  +3.83%  " }
    15 = @{ Price = "3.297.07"; Volume = "  +3.25%  " }
    16 = @{ Price = "2.848.58"; Volume = "  +2.93%  " }
    17 = @{ Price = "0.908"; Volume = "  +2.67%  " }
    18 = @{ Price = "52.660.18"; Volume = "  +1.76%  " }
    19 = @{ Price = "3.21"; Volume = "  +7.12%  " }
    20 = @{ Price = "7.34"; Volume = "  +7.03%  " }
    21 = @{ Price = "13.83"; Volume = "  +1.42%  " }
    22 = @{ Price = "0.0₃0990"; Volume = "  +2.73%  " }
    23 = @{ Price = "273.03"; Volume = "  -2.16%  " }
    24 = @{ Price = "70.79"; Volume = "  +1.34%  " }
    25 = @{ Price = "2.85"; Volume = "  +7.92%  " }
    26 = @{ Price = "27.30"; Volume = "  +1.63%  " }
    27 = @{ Volume = "  +0.01%  " }
    28 = @{ Price = "10.36"; Volume = "  +1.02%  " }
    29 = @{ Volume = "  +1.70%  " }
    30 = @{ Volume = "  +0.83%  " }
    31 = @{ Price = "34.84"; Volume = "  -1.30%  " }
    32 = @{ Price = "51.10"; Volume = "  +1.29%  " }
    33 = @{ Price = "0.0454"; Volume = "  +31.48%  " }
    34 = @{ Price = "5.86"; Volume = "  +4.81%  " }
    35 = @{ Price = "0.0839"; Volume = "  +1.99%  " }
    36 = @{ Volume = "  +1.59%  " }
    37 = @{ Price = "1.00"; Volume = "  +0.09%  " }
    38 = @{ Price = "18.92"; Volume = "  -1.24%  " }
    39 = @{ Price = "3.30"; Volume = "  +2.40%  " }
    40 = @{ Price = "5.00"; Volume = "  -0.08%  " }
    41 = @{ Volume = "  +10.90%  " }
    42 = @{ Price = "23.90"; Volume = "  +2.82%  " }
    43 = @{ Price = "0.117"; Volume = "  +2.85%  " }
    44 = @{ Price = "126.68"; Volume = "  -3.76%  " }
    45 = @{ Price = "2.30"; Volume = "  +0.65%  " }
    46 = @{ Price = "3.39"; Volume = "  +1.33%  " }
    47 = @{ Price = "2.084.41"; Volume = "  -1.24%  " }
    48 = @{ Volume = "  +2.65%  " }
    49 = @{ Price = "0.969"; Volume = "  +10.68%  " }
    50 = @{ Price = "5.67"; Volume = "  +1.21%  " }
    51 = @{ Price = "9.01"; Volume = "  +0.26%  " }
}

foreach ($row in $updates.Keys) {
    $entry = $updates[$row]
    if ($entry.Contains("Price")) {
        Set-TextValue "D$row" $entry.Price
    }
    $ws.Range("E$row").Value = $entry.Volume
}
